# Generate Report for Handoff
# Updates the localization-status workbook with a fresh handoff generation:
#  - file name guid changes from 080c98f6-... to 43676303-...
#  - xliff hash changes from 384be33... to d039b541...
#  - handoff timestamps move forward; handback data is cleared (new handoff, no handback yet)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "43676303-30b4-4831-90b3-9c90e667589b.md"
$overview.Range("B2").Value = "e2e\43676303-30b4-4831-90b3-9c90e667589b.md"
$overview.Range("G2").Value = "2016-08-17 16:58:36"

$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c16cba3554c1087eb3d4c7e33db5ad31c17fb5e7/e2e/43676303-30b4-4831-90b3-9c90e667589b.md", "", "", "e2e\43676303-30b4-4831-90b3-9c90e667589b.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("G2").Value = "43676303-30b4-4831-90b3-9c90e667589b.d039b5418b4a2a03251de53753f241e25e77ae84.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-17 16:58:31"
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

# Latest Target File / Latest Handback File are reset (no handback has happened yet)
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c16cba3554c1087eb3d4c7e33db5ad31c17fb5e7/e2e/43676303-30b4-4831-90b3-9c90e667589b.md", "", "", "43676303-30b4-4831-90b3-9c90e667589b.md")

$zhcn.Range("I2").Style = "Normal"
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""

$zhcn.Columns.Item(9).ColumnWidth = 18.6506053379604
$zhcn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("G2").Value = "43676303-30b4-4831-90b3-9c90e667589b.d039b5418b4a2a03251de53753f241e25e77ae84.de-de.xlf"
$dede.Range("K2").Value = "0001-01-01 00:00:00"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c16cba3554c1087eb3d4c7e33db5ad31c17fb5e7/e2e/43676303-30b4-4831-90b3-9c90e667589b.md", "", "", "43676303-30b4-4831-90b3-9c90e667589b.md")

$dede.Range("I2").Style = "Normal"
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""

$dede.Columns.Item(9).ColumnWidth = 18.6506053379604
$dede.Columns.Item(10).ColumnWidth = 21.7054770333426
